# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect a newer data pull, per the diff.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F5").Value = 436
$sheet1.Range("F6").Value = 81
$sheet1.Range("F7").Value = 564
$sheet1.Range("F8").Value = 75
$sheet1.Range("F9").Value = 6833
$sheet1.Range("F10").Value = 161
$sheet1.Range("F11").Value = 99
$sheet1.Range("F15").Value = 1107
$sheet1.Range("F16").Value = 16254
$sheet1.Range("F17").Value = 1598
$sheet1.Range("F18").Value = 44
$sheet1.Range("F22").Value = 11397
$sheet1.Range("F24").Value = 1042
$sheet1.Range("F25").Value = 4486
$sheet1.Range("F26").Value = 333
$sheet1.Range("F28").Value = 52

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F5").Value = 436
$sheet4.Range("F6").Value = 81
$sheet4.Range("F7").Value = 564
$sheet4.Range("F9").Value = 75
$sheet4.Range("F10").Value = 6833
$sheet4.Range("F11").Value = 161
$sheet4.Range("F12").Value = 99
$sheet4.Range("F17").Value = 1107
$sheet4.Range("F18").Value = 16254
$sheet4.Range("F19").Value = 1598
$sheet4.Range("F20").Value = 44
$sheet4.Range("F26").Value = 11397
$sheet4.Range("F28").Value = 1042
$sheet4.Range("F29").Value = 4486
$sheet4.Range("F30").Value = 333
$sheet4.Range("F32").Value = 52
